{"js": "// 1) Append a new sentence to the end of the \"Tour\" paragraph (same paragraph,\n//    new run), right after \"... beschikbaarheid van kaartjes ziet.\"\nconst tourMatches = context.document.body.search(\"beschikbaarheid van kaartjes ziet.\", { matchCase: false });\nawait context.sync();\n\nif (tourMatches.items.length > 0) {\n  const tourParagraph = tourMatches.items[0].paragraphs.getFirst();\n  tourParagraph.insertText(\n    \" Je kan ook op More Information klikken waardoor je naar een pagina gaat waar meer details opstaat.\",\n    \"End\"\n  );\n  await context.sync();\n}\n\n// 2) Add a new sentence right before the \"_GoBack\" bookmark and another right\n//    after it, keeping the bookmark itself in place (in the middle of the\n//    paragraph). Inserting \"After\" first \u2014 while the bookmark is still at its\n//    original (untouched) location \u2014 then \"Before\" keeps the bookmark anchored\n//    between the two new runs instead of being pushed to one end.\nconst bookmarkName = \"_GoBack\";\nconst afterRange = context.document.getBookmarkRange(bookmarkName);\nafterRange.insertText(\" en verwijderen\", \"After\");\nawait context.sync();\n\nconst beforeRange = context.document.getBookmarkRange(bookmarkName);\nbeforeRange.insertText(\n  \"Op de CMS pagina, beschikbaar door op Alt + L te drukken, kan je content toevoegen, aanpassen\",\n  \"Before\"\n);\nawait context.sync();\n\n// 3) Append a blank paragraph, a paragraph about pagination, and another\n//    blank paragraph at the very end of the document body.\nconst body = context.document.body;\nbody.insertParagraph(\"\", \"End\");\nawait context.sync();\n\nbody.insertParagraph(\"Op de pagina\\u2019s waar het nodig is, is gebruik gemaakt van pagination.\", \"End\");\nawait context.sync();\n\nbody.insertParagraph(\"\", \"End\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Append a new sentence to the end of the \"Tour\" paragraph (same\n#    paragraph, new run), right after \"... beschikbaarheid van kaartjes ziet.\"\n$tourRange = $d.Content.Find.Execute(\"beschikbaarheid van kaartjes ziet.\")\n$found = $d.Content\n$found.Find.Text = \"beschikbaarheid van kaartjes ziet.\"\nif ($found.Find.Execute()) {\n    $tourParagraph = $found.Paragraphs.Item(1)\n    $tourParagraph.Range.InsertAfter(\" Je kan ook op More Information klikken waardoor je naar een pagina gaat waar meer details opstaat.\")\n}\n\n# 2) Add a new sentence right before the \"_GoBack\" bookmark and another right\n#    after it, keeping the bookmark itself in place (in the middle of the\n#    paragraph). Inserting \"After\" first -- while the bookmark is still at its\n#    original (untouched) location -- then \"Before\" keeps the bookmark\n#    anchored between the two new runs instead of being pushed to one end.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Range.InsertAfter(\" en verwijderen\")\n\n$bm2 = $d.Bookmarks.Item(\"_GoBack\")\n$bm2.Range.InsertBefore(\"Op de CMS pagina, beschikbaar door op Alt + L te drukken, kan je content toevoegen, aanpassen\")\n\n# 3) Append a blank paragraph, a paragraph about pagination, and another\n#    blank paragraph at the very end of the document body.\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n\n$paginationParagraph = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n$paginationParagraph.Range.InsertBefore(\"Op de pagina\u2019s waar het nodig is, is gebruik gemaakt van pagination.\")\n"}
